$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.049772999999999
$ws.Range("H2").Value = 15.149319
$ws.Range("I2").Value = 0.1400646900514762
$ws.Range("J2").Value = 0.1400646900514762
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.325274
$ws.Range("N2").Value = 0.975822
$ws.Range("O2").Value = 0.4689956999283401
$ws.Range("P2").Value = 0.4689956999283402
$ws.Range("Q2").Value = 1.642559862802
$ws.Range("R2").Value = 14.783038765218
$ws.Range("S2").Value = 0.0656897373459381
$ws.Range("T2").Value = 0.06568973734593811

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.049772999999999
$ws.Range("H3").Value = 15.149319
$ws.Range("I3").Value = 0.1400646900514762
$ws.Range("J3").Value = 0.1400646900514762
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3682803333333333
$ws.Range("N3").Value = 1.104841
$ws.Range("O3").Value = 0.5310043000716598
$ws.Range("P3").Value = 0.5310043000716599
$ws.Range("Q3").Value = 1.859732083697666
$ws.Range("R3").Value = 16.737588753279
$ws.Range("S3").Value = 0.0743749527055381
$ws.Range("T3").Value = 0.07437495270553811

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.420946
$ws.Range("H4").Value = 61.262838
$ws.Range("I4").Value = 0.5664122866607931
$ws.Range("J4").Value = 0.5664122866607931
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.325274
$ws.Range("N4").Value = 0.975822
$ws.Range("O4").Value = 0.4689956999283401
$ws.Range("P4").Value = 0.4689956999283402
$ws.Range("Q4").Value = 6.642402789204001
$ws.Range("R4").Value = 59.781625102836
$ws.Range("S4").Value = 0.2656449268304903
$ws.Range("T4").Value = 0.2656449268304903

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 20.420946
$ws.Range("H5").Value = 61.262838
$ws.Range("I5").Value = 0.5664122866607931
$ws.Range("J5").Value = 0.5664122866607931
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3682803333333333
$ws.Range("N5").Value = 1.104841
$ws.Range("O5").Value = 0.5310043000716598
$ws.Range("P5").Value = 0.5310043000716599
$ws.Range("Q5").Value = 7.520632799862
$ws.Range("R5").Value = 67.685695198758
$ws.Range("S5").Value = 0.3007673598303028
$ws.Range("T5").Value = 0.3007673598303028

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.258872
$ws.Range("H6").Value = 30.776616
$ws.Range("I6").Value = 0.2845485781158417
$ws.Range("J6").Value = 0.2845485781158416
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.325274
$ws.Range("N6").Value = 0.975822
$ws.Range("O6").Value = 0.4689956999283401
$ws.Range("P6").Value = 0.4689956999283402
$ws.Range("Q6").Value = 3.336944330928
$ws.Range("R6").Value = 30.032498978352
$ws.Range("S6").Value = 0.1334520595570531
$ws.Range("T6").Value = 0.1334520595570531

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.258872
$ws.Range("H7").Value = 30.776616
$ws.Range("I7").Value = 0.2845485781158417
$ws.Range("J7").Value = 0.2845485781158416
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3682803333333333
$ws.Range("N7").Value = 1.104841
$ws.Range("O7").Value = 0.5310043000716598
$ws.Range("P7").Value = 0.5310043000716599
$ws.Range("Q7").Value = 3.778140799784
$ws.Range("R7").Value = 34.003267198056
$ws.Range("S7").Value = 0.1510965185587885
$ws.Range("T7").Value = 0.1510965185587885

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.323557
$ws.Range("H8").Value = 0.9706710000000001
$ws.Range("I8").Value = 0.008974445171889013
$ws.Range("J8").Value = 0.008974445171889013
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.325274
$ws.Range("N8").Value = 0.975822
$ws.Range("O8").Value = 0.4689956999283401
$ws.Range("P8").Value = 0.4689956999283402
$ws.Range("Q8").Value = 0.105244679618
$ws.Range("R8").Value = 0.947202116562
$ws.Range("S8").Value = 0.0042089761948586
$ws.Range("T8").Value = 0.0042089761948586

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.323557
$ws.Range("H9").Value = 0.9706710000000001
$ws.Range("I9").Value = 0.008974445171889013
$ws.Range("J9").Value = 0.008974445171889013
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3682803333333333
$ws.Range("N9").Value = 1.104841
$ws.Range("O9").Value = 0.5310043000716598
$ws.Range("P9").Value = 0.5310043000716599
$ws.Range("Q9").Value = 0.1191596798123333
$ws.Range("R9").Value = 1.072437118311
$ws.Range("S9").Value = 0.004765468977030412
$ws.Range("T9").Value = 0.004765468977030412

